$d = $word.ActiveDocument

# --- Table of functionalities: add two new rows -------------------------
# Table #2 in the document ("Funcionalidad" / "Estado de implementación")
$t = $d.Tables.Item(2)

# Insert "Ver historial de compra" / "0%" right before the existing last
# row ("Modificar el inventario, quitar o agregar productos (administrador)").
$lastRow = $t.Rows.Item($t.Rows.Count)
$newRow1 = $t.Rows.Add($lastRow)
$newRow1.Cells.Item(1).Range.Text = "Ver historial de compra"
$newRow1.Cells.Item(2).Range.Text = "0%"

# Append "Acceder a los datos de todos los clientes registrados
# (administrador)" / "10%" as the new final row of the table.
$newRow2 = $t.Rows.Add()
$newRow2.Cells.Item(1).Range.Text = "Acceder a los datos de todos los clientes registrados (administrador)"
$newRow2.Cells.Item(2).Range.Text = "10%"

# --- Remove the stray _GoBack bookmark near the end of the document -----
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
